# Coastal Surface Piercing Profilers - update Coastal CSPP ingest/cal sheets
# Change reference designator from GI05MOAS-GL001 to GI05MOAS-GL469
# across the "Moorings" and "Asset_Cal_Info" sheets, and update the
# active selection on the "Moorings" sheet.

$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsCalInfo  = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet: update Ref Des for deployment row ---
$wsMoorings.Range("A2").Value = "GI05MOAS-GL469"

# --- Asset_Cal_Info sheet: update Ref Des / instrument reference designators ---
$wsCalInfo.Range("A2").Value  = "GI05MOAS-GL469-01-FLORDM000"
$wsCalInfo.Range("A3").Value  = "GI05MOAS-GL469-01-FLORDM000"
$wsCalInfo.Range("A4").Value  = "GI05MOAS-GL469-01-FLORDM000"
$wsCalInfo.Range("A5").Value  = "GI05MOAS-GL469-01-FLORDM000"
$wsCalInfo.Range("A7").Value  = "GI05MOAS-GL469-02-DOSTAM000"
$wsCalInfo.Range("A9").Value  = "GI05MOAS-GL469-04-CTDGVM000"
$wsCalInfo.Range("A11").Value = "GI05MOAS-GL469-00-ENG000000"

# --- Moorings sheet: move active selection from B2 to D30 ---
$wsMoorings.Activate()
$wsMoorings.Range("D30").Select()
